# Reorders the species-record data held in rows 67, 68, 70, 71, 72 and 73
# of the "Artfynd" sheet. Row 69 is left untouched. The row-number/location
# related columns (C, I, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX,
# AY) are identical across all six rows, so only A, B, D, E, F, G, H, K, M,
# P, Q and R need to move.
#
# The permutation (new row <- source row):
#   67 <- 72      (67 and 72 simply swap)
#   72 <- 67
#   68 <- 73      (68, 73, 70, 71 rotate in a 4-cycle)
#   73 <- 70
#   70 <- 71
#   71 <- 68

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Snapshot the "before" values of every cell we are going to move,
#    so that later writes don't clobber data we still need to read.
# ---------------------------------------------------------------------
$cols = @("A","B","D","E","F","G","H","P","Q","R")
$rows = @(67,68,70,71,72,73)

$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value()
    }
    $snapshot[$r] = $rowData
}

# ---------------------------------------------------------------------
# 2. Write the permuted values back.
# ---------------------------------------------------------------------
$mapping = @{ 67 = 72; 68 = 73; 70 = 71; 71 = 68; 72 = 67; 73 = 70 }

foreach ($target in $rows) {
    $source = $mapping[$target]
    $data = $snapshot[$source]
    foreach ($c in $cols) {
        $ws.Range("$c$target").Value = $data[$c]
    }
}

# ---------------------------------------------------------------------
# 3. Column K ("Ålder-Stadium") - present-but-empty marker cell.
#      67: absent  -> blank-present
#      68: absent  -> blank-present
#      70: blank-present -> blank-present (no change needed)
#      71: blank-present -> absent
#      72: blank-present -> absent
#      73: blank-present -> blank-present (no change needed)
# ---------------------------------------------------------------------
$ws.Range("K70").Copy($ws.Range("K67"))
$ws.Range("K70").Copy($ws.Range("K68"))
$ws.Range("K71").ClearContents()
$ws.Range("K72").ClearContents()

# ---------------------------------------------------------------------
# 4. Column M ("Aktivitet") - free text, only present on two rows.
#      68: "gammalt bo" -> absent
#      71: absent -> "gammalt bo"
# ---------------------------------------------------------------------
$ws.Range("M68").ClearContents()
$ws.Range("M71").Value = "gammalt bo"
